$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New checklist rows appended below the existing last row (15), matching
# the formatting of row 15 (Good/"passed" style, 34pt row height).
$newRows = @(
    @{ Row = 16; Num = 1.15; Name = "Проверка появления проекта в списке проектов после его создания" },
    @{ Row = 17; Num = 1.16; Name = "Удаление проекта из списка проектов после нажатия кнопки удалить" },
    @{ Row = 18; Num = 1.17; Name = "Проверка правильного отображения задач прикрепленных к этому проекту" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Copy formatting from the row above (cells only, not whole rows, so we
    # don't stamp formats across the unused columns) and set new values.
    $ws.Range("A" + ($row - 1) + ":C" + ($row - 1)).Copy()
    $ws.Range("A" + $row + ":C" + $row).PasteSpecial(-4122) # xlPasteFormats
    $ws.Rows.Item($row).RowHeight = $ws.Rows.Item($row - 1).RowHeight

    $ws.Cells.Item($row, 1).Value = $r.Num
    $ws.Cells.Item($row, 2).Value = $r.Name
    $ws.Cells.Item($row, 3).Value = "passed"
}

$excel.CutCopyMode = 0
$ws.Range("B18").Select()
